$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Move the small price-legend block (currently C15:D16) up to E7:F8 so the
#    rows it used to occupy are free for the new budget lines.
# ---------------------------------------------------------------------------
$ws.Range("E7").Value = $ws.Range("C15").Value2
$ws.Range("F7").Value = $ws.Range("D15").Value2
$ws.Range("E8").Value = $ws.Range("C16").Value2
$ws.Range("F8").Value = $ws.Range("D16").Value2

$ws.Range("E7:E8").HorizontalAlignment = -4108   # xlCenter
$ws.Range("F7:F8").HorizontalAlignment = -4108   # xlCenter
$ws.Range("F7:F8").VerticalAlignment = -4108     # xlCenter

# Remove the now-unused old legend cells entirely (content + formatting).
$ws.Range("C15:D16").Clear()

$ws.Columns.Item("E").ColumnWidth = 31.7109375
$ws.Columns.Item("F").ColumnWidth = 13.42578125

# ---------------------------------------------------------------------------
# 2. Header row: center the two titles.
# ---------------------------------------------------------------------------
$ws.Range("C1").HorizontalAlignment = -4108      # splits off a new style for C1
$ws.Range("B1").HorizontalAlignment = -4108      # only B1 used that style -> updates it directly

# ---------------------------------------------------------------------------
# 3. Move the TOTAL row from row 13 down to row 17, keeping its look, then
#    recompute the sum range.
# ---------------------------------------------------------------------------
$ws.Range("B13").HorizontalAlignment = -4108     # only B13 used that style -> updates it directly

$ws.Range("C1").Copy()
$ws.Range("C17").PasteSpecial(-4122)             # xlPasteFormats
$ws.Range("B13").Copy()
$ws.Range("B17").PasteSpecial(-4122)             # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B17").Value = "TOTAL"
$ws.Range("C17").Formula = "=SUM(C2:C16)"
$ws.Rows.Item(17).RowHeight = 15.75

$ws.Range("B13:C13").Clear()

# ---------------------------------------------------------------------------
# 4. Fill the new budget rows 13-16 *before* touching row 12's border, so row
#    16 can still borrow row 12's current (thick-bottom) look. Row 16's text
#    is entered before row 15's so the shared-string table ends up ordered
#    the same way Excel saved it.
# ---------------------------------------------------------------------------
$ws.Range("B9:C9").Copy()
$ws.Range("B13:C15").PasteSpecial(-4122)         # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B12:C12").Copy()
$ws.Range("B16:C16").PasteSpecial(-4122)         # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B13").Value = "Brief Charte graphique (graphiste)"
$ws.Range("C13").Formula = "=900"

$ws.Range("B14").Value = "Formation des commerçants (equipe de 5 personnes)"
$ws.Range("C14").Formula = "=5*500"

$ws.Range("B16").Value = "Référencement (dev)"
$ws.Range("C16").Value = 500
$ws.Rows.Item(16).RowHeight = 15.75

$ws.Range("B15").Value = "Publicité sur les réseaux (community manager)"
$ws.Range("C15").Value = 500

# ---------------------------------------------------------------------------
# 5. "Mise en ligne" (row 12) loses its thick bottom border - it becomes a
#    plain row like the others above it.
# ---------------------------------------------------------------------------
$ws.Range("B9:C9").Copy()
$ws.Range("B12:C12").PasteSpecial(-4122)         # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows.Item(12).AutoFit()

# ---------------------------------------------------------------------------
# 6. Selection marker, matching what Excel saved last.
# ---------------------------------------------------------------------------
$ws.Range("E16").Select()
